# "Add files via upload" — three new keyword-stem rows were added to the
# Cluster_Keywords table (Resta -> Food & Drink; Conce & Squar -> Retail),
# the table was re-sorted by Cluster Category then Stem (its existing sort
# order), and the selection moved to the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

function Add-ClusterRow {
    param($Stem, $Category)

    $newRow = $lo.ListRows.Add()
    $r = $newRow.Range.Row

    $ws.Cells.Item($r, 1).Value = $Stem
    $ws.Cells.Item($r, 2).Formula = "=LEN(Cluster_Keywords[[#This Row],[Stem]])"
    $ws.Cells.Item($r, 3).Value = $Category

    # Pick up the same cell style ("s=1") the rest of the data rows use.
    $srcFormat = $ws.Range("A2:C2")
    $dstFormat = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 3))
    $srcFormat.Copy()
    $dstFormat.PasteSpecial(-4122)
}

Add-ClusterRow "Resta" "Food & Drink"
Add-ClusterRow "Conce" "Retail"
Add-ClusterRow "Squar" "Retail"

# Re-apply the table's existing sort (Cluster Category, then Stem) now that
# the new rows have been appended at the bottom.
$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($lo.ListColumns.Item("Cluster Category").Range)
$lo.Sort.SortFields.Add($lo.ListColumns.Item("Stem").Range)
$lo.Sort.Apply()

# The sort can leave the calculated column referring to [@Stem]; restate it
# explicitly so every row keeps the table's own calculated-column formula.
for ($r = 2; $r -le $lo.Range.Rows.Count; $r++) {
    $ws.Cells.Item($r, 2).Formula = "=LEN(Cluster_Keywords[[#This Row],[Stem]])"
}

# Make sure every data row (including the three new ones) still has the
# "highlight duplicate values" conditional formatting the rest of column A
# carries, in the same red-on-pink style as the existing rule.
$dataCol = $ws.Range("A2:A" + $lo.Range.Rows.Count)
$dupRule = $dataCol.FormatConditions.AddUniqueValues()
$dupRule.DupeUnique = 1
$dupRule.Font.Color = 393372
$dupRule.Interior.Color = 13551615
$dupRule.SetFirstPriority()

$ws.Range("A" + $lo.Range.Rows.Count).Select()
